$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: price per unit + total (quantity stays 1)
$ws.Range("B2").Value = 4869
$ws.Range("D2").Value = 4869

# Row 3: price per unit + total (quantity stays 12)
$ws.Range("B3").Value = 3517
$ws.Range("D3").Value = 42204

# Row 4: price per unit + total (quantity stays 1)
$ws.Range("B4").Value = 4328
$ws.Range("D4").Value = 4328
